# Update dispensing data from FY2022/23 to FY2023/24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the fy (year) column
$ws.Range("A2").Value = "2023/24"
$ws.Range("A3").Value = "2023/24"

# Update South West row (row 2): pop, n_pharmacies, n_items
$ws.Range("C2").Value = 5811259
$ws.Range("D2").Value = 1107
$ws.Range("E2").Value = 99164702

# Update England row (row 3): pop, n_pharmacies, n_items
$ws.Range("C3").Value = 57690323
$ws.Range("D3").Value = 12009
$ws.Range("E3").Value = 1112920890

# Update the selected cell (matches author's saved selection state)
$ws.Range("C3").Select()
